$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02902903341425
$ws.Range("D2").Value = 1.032114794810221
$ws.Range("E2").Value = 1.02887243402806
$ws.Range("I2").Value = 1.030244821884919
$ws.Range("J2").Value = 1.034178074952713
$ws.Range("K2").Value = 1.034920958377131
$ws.Range("L2").Value = 1.031687988137549
$ws.Range("N2").Value = 1.035646726441037
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030499042817314
$ws.Range("D3").Value = 1.033207224475991
$ws.Range("E3").Value = 1.030136745897304
$ws.Range("I3").Value = 1.030523109921155
$ws.Range("J3").Value = 1.035286280431448
$ws.Range("K3").Value = 1.03582133351086
$ws.Range("L3").Value = 1.032759097327895
$ws.Range("N3").Value = 1.036756505698657
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03144863523407
$ws.Range("D4").Value = 1.033912465228322
$ws.Range("E4").Value = 1.030953718557134
$ws.Range("I4").Value = 1.030701133404464
$ws.Range("J4").Value = 1.036001403749696
$ws.Range("K4").Value = 1.036401720575996
$ws.Range("L4").Value = 1.03345051212698
$ws.Range("N4").Value = 1.037472644574042
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031847467667991
$ws.Range("D5").Value = 1.034208561536003
$ws.Range("E5").Value = 1.031296910621233
$ws.Range("I5").Value = 1.030775485798795
$ws.Range("J5").Value = 1.036301577599242
$ws.Range("K5").Value = 1.036645188890909
$ws.Range("L5").Value = 1.033740789260599
$ws.Range("N5").Value = 1.037773244704885
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031914411538146
$ws.Range("D6").Value = 1.034258254846368
$ws.Range("E6").Value = 1.031354518795202
$ws.Range("I6").Value = 1.030787941279796
$ws.Range("J6").Value = 1.036351951050152
$ws.Range("K6").Value = 1.036686037541243
$ws.Range("L6").Value = 1.033789505124613
$ws.Range("N6").Value = 1.037823689691874
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031453965921846
$ws.Range("D7").Value = 1.033916423196602
$ws.Range("E7").Value = 1.030958305334133
$ws.Range("I7").Value = 1.030702128823139
$ws.Range("J7").Value = 1.036005416504179
$ws.Range("K7").Value = 1.036404975874957
$ws.Range("L7").Value = 1.033454392365066
$ws.Range("N7").Value = 1.037476663027097
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029526166075662
$ws.Range("D8").Value = 1.032484327548239
$ws.Range("E8").Value = 1.02929994958345
$ws.Range("I8").Value = 1.030339295478387
$ws.Range("J8").Value = 1.034553007080374
$ws.Range("K8").Value = 1.035225705636768
$ws.Range("L8").Value = 1.032050322059793
$ws.Range("N8").Value = 1.036022191015317
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026116539187799
$ws.Range("D9").Value = 1.029948079159845
$ws.Range("E9").Value = 1.026368886204631
$ws.Range("I9").Value = 1.029684187883562
$ws.Range("J9").Value = 1.031978434011029
$ws.Range("K9").Value = 1.033130526297976
$ws.Range("L9").Value = 1.029563205358117
$ws.Range("N9").Value = 1.033443961756906
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023834523652523
$ws.Range("D10").Value = 1.028248422099237
$ws.Range("E10").Value = 1.02440857245048
$ws.Range("I10").Value = 1.029236764775878
$ws.Range("J10").Value = 1.03025147569781
$ws.Range("K10").Value = 1.031721946703666
$ws.Range("L10").Value = 1.027896103553927
$ws.Range("N10").Value = 1.031714550964797
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022844164068746
$ws.Range("D11").Value = 1.027510296056301
$ws.Range("E11").Value = 1.023558171951071
$ws.Range("I11").Value = 1.029040469508607
$ws.Range("J11").Value = 1.029501101268163
$ws.Range("K11").Value = 1.031109160350854
$ws.Range("L11").Value = 1.027172020668955
$ws.Range("N11").Value = 1.030963110917391
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022475956125085
$ws.Range("D12").Value = 1.027235792658514
$ws.Range("E12").Value = 1.023242052613344
$ws.Range("I12").Value = 1.028967170427715
$ws.Range("J12").Value = 1.029221983199112
$ws.Range("K12").Value = 1.030881109063103
$ws.Range("L12").Value = 1.026902724955423
$ws.Range("N12").Value = 1.030683596468667
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022554953704357
$ws.Range("D13").Value = 1.027294689614867
$ws.Range("E13").Value = 1.023309872374868
$ws.Range("I13").Value = 1.028982910840851
$ws.Range("J13").Value = 1.029281872978336
$ws.Range("K13").Value = 1.030930046613285
$ws.Range("L13").Value = 1.026960505242221
$ws.Range("N13").Value = 1.030743571298246
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022813734961144
$ws.Range("D14").Value = 1.027487612294362
$ws.Range("E14").Value = 1.023532046417946
$ws.Range("I14").Value = 1.029034418475041
$ws.Range("J14").Value = 1.02947803738913
$ws.Range("K14").Value = 1.031090318467228
$ws.Range("L14").Value = 1.027149767553738
$ws.Range("N14").Value = 1.030940014285004
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02297313274105
$ws.Range("D15").Value = 1.027606434361677
$ws.Range("E15").Value = 1.023668902816389
$ws.Range("I15").Value = 1.029066102784046
$ws.Range("J15").Value = 1.02959884814226
$ws.Range("K15").Value = 1.031189009429178
$ws.Range("L15").Value = 1.027266333195938
$ws.Range("N15").Value = 1.03106099660326
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023900202768807
$ws.Range("D16").Value = 1.028297363046313
$ws.Range("E16").Value = 1.024464977055054
$ws.Range("I16").Value = 1.029249738183741
$ws.Range("J16").Value = 1.030301220456818
$ws.Range("K16").Value = 1.031762554575832
$ws.Range("L16").Value = 1.02794411133041
$ws.Range("N16").Value = 1.031764366367069
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024481125460074
$ws.Range("D17").Value = 1.028730181302901
$ws.Range("E17").Value = 1.024963907962537
$ws.Range("I17").Value = 1.029364241568274
$ws.Range("J17").Value = 1.030741102067856
$ws.Range("K17").Value = 1.032121554405474
$ws.Range("L17").Value = 1.028368665898162
$ws.Range("N17").Value = 1.032204872660449
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024819753113315
$ws.Range("D18").Value = 1.028982428515723
$ws.Range("E18").Value = 1.025254774721903
$ws.Range("I18").Value = 1.029430782768026
$ws.Range("J18").Value = 1.030997428104095
$ws.Range("K18").Value = 1.032330677133191
$ws.Range("L18").Value = 1.028616087916753
$ws.Range("N18").Value = 1.032461562709061
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024935180241989
$ws.Range("D19").Value = 1.02906840312172
$ws.Range("E19").Value = 1.025353927322819
$ws.Range("I19").Value = 1.029453429813726
$ws.Range("J19").Value = 1.031084786539189
$ws.Range("K19").Value = 1.032401935953746
$ws.Range("L19").Value = 1.028700416403255
$ws.Range("N19").Value = 1.032549045203152
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024418820228211
$ws.Range("D20").Value = 1.028683765591045
$ws.Range("E20").Value = 1.024910393084926
$ws.Range("I20").Value = 1.029351981968938
$ws.Range("J20").Value = 1.030693932784218
$ws.Range("K20").Value = 1.03208306568587
$ws.Range("L20").Value = 1.028323137335802
$ws.Range("N20").Value = 1.032157636391017
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022737539907905
$ws.Range("D21").Value = 1.027430810542181
$ws.Range("E21").Value = 1.023466628453551
$ws.Range("I21").Value = 1.029019261439518
$ws.Range("J21").Value = 1.029420282831989
$ws.Range("K21").Value = 1.031043134454876
$ws.Range("L21").Value = 1.027094043929562
$ws.Range("N21").Value = 1.030882177709767
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021678455727081
$ws.Range("D22").Value = 1.026641112918962
$ws.Range("E22").Value = 1.022557468667462
$ws.Range("I22").Value = 1.028807830967081
$ws.Range("J22").Value = 1.028617196519454
$ws.Range("K22").Value = 1.030386768171456
$ws.Range("L22").Value = 1.026319299220231
$ws.Range("N22").Value = 1.030077950922553
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022240088557905
$ws.Range("D23").Value = 1.027059929946879
$ws.Range("E23").Value = 1.02303956711845
$ws.Range("I23").Value = 1.028920126852059
$ws.Range("J23").Value = 1.029043147199596
$ws.Range("K23").Value = 1.030734960980156
$ws.Range("L23").Value = 1.026730194329263
$ws.Range("N23").Value = 1.030504506501519
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024446973944305
$ws.Range("D24").Value = 1.02870473949737
$ws.Range("E24").Value = 1.024934574622061
$ws.Range("I24").Value = 1.029357522316518
$ws.Range("J24").Value = 1.030715247324177
$ws.Range("K24").Value = 1.032100457933955
$ws.Range("L24").Value = 1.02834371039324
$ws.Range("N24").Value = 1.032178981200069
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026999549576349
$ws.Range("D25").Value = 1.030605295096803
$ws.Range("E25").Value = 1.027127716964246
$ws.Range("I25").Value = 1.029855425047492
$ws.Range("J25").Value = 1.032645862757888
$ws.Range("K25").Value = 1.03367423985535
$ws.Range("L25").Value = 1.030207752532584
$ws.Range("N25").Value = 1.034112338329141
